# Update task assignment (Resource) and Deadline columns on Sheet2,
# and add a new (empty) Sheet3 at the end of the workbook.

$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("Sheet2")

# Insert a new row 1, pushing the existing header/data rows down by one.
$ws2.Rows.Item(1).Insert()

# Resource assignments (column D) and Deadlines (column E), grouped by
# contiguous ranges that share the same values.
$groups = @(
    @(2, 6, "Trung", "6/4/2014 12:00AM"),
    @(7, 11, "Nhan", "6/4/2014 12:00AM"),
    @(12, 15, "Thang", "6/4/2014 12:00AM"),
    @(16, 21, "Hung", "6/4/2014 12:00AM"),
    @(22, 34, "Nhan", "8/4/2013 12:00AM"),
    @(35, 36, "Thang", "8/4/2013 12:00AM"),
    @(37, 39, "Trung", "8/4/2013 12:00AM"),
    @(40, 47, "Hung", "8/4/2013 12:00AM"),
    @(48, 48, "Thang", "8/4/2013 12:00AM")
)

# Fill column D (Resource) first, in row order, so the new shared-string
# table picks up "Trung"/"Nhan"/"Thang"/"Hung" in first-seen order.
foreach ($g in $groups) {
    $startRow = $g[0]
    $endRow = $g[1]
    $resource = $g[2]

    $ws2.Range("D$startRow`:D$endRow").Value = $resource
}

# New header row for the Resource / Deadline columns.
$ws2.Range("D1").Value = "Resource"
$ws2.Range("E1").Value = "Deadline"

# Now fill column E (Deadline) in row order.
foreach ($g in $groups) {
    $startRow = $g[0]
    $endRow = $g[1]
    $deadline = $g[3]

    $ws2.Range("E$startRow`:E$endRow").Value = $deadline
}

# Trailing blank rows (49-53) only carry the deadline-column date style.
$ws2.Range("E49:E53").NumberFormat = "d-mmm"

# Apply the date number format to the whole Deadline column range used above.
$ws2.Range("E2:E48").NumberFormat = "d-mmm"

# Size the new Deadline column to fit its (text) contents.
$ws2.Columns.Item(5).AutoFit()

# Add a new, empty worksheet ("Sheet3") at the end of the workbook.
$lastIndex = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($lastIndex)
$ws3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws3.Name = "Sheet3"

# Sheet2 stays the selected/active tab; reflect the recorded scroll
# position and selection.
$ws2.Activate()
$ws2.Application.Goto($ws2.Range("A32"), $true)
$ws2.Range("F47").Select()
